$d = $word.ActiveDocument

# --- 1) "Name" run: remove the bold direct formatting -----------------
$rName = $d.Content
$rName.Find.Execute("Name")
$rName.Font.Bold = 0

# --- 2) "Address" run: remove the underline direct formatting ---------
$rAddr = $d.Content
$rAddr.Find.Execute("Address")
$rAddr.Font.Underline = 0

# --- 3) Merge the two runs that spell out ": <# <C" + "ontent " -------
#     into a single run ": <# <Content " (also absorbs/repositions the
#     "_GoBack" bookmark that previously sat between them).
$rTail = $d.Content
$rTail.Find.Execute(": <# <Content ")
$rTail.Text = "TEMP_MERGE_PLACEHOLDER"

$rTail2 = $d.Content
$rTail2.Find.Execute("TEMP_MERGE_PLACEHOLDER")
$rTail2.Text = ": <# <Content "

# --- 4) Re-anchor the "_GoBack" bookmark around the word "Address" ----
$rAddr2 = $d.Content
$rAddr2.Find.Execute("Address")
$d.Bookmarks.Add("_GoBack", $rAddr2)
